$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The author re-exported this results table after dropping a couple of
# now-unused parameters from rdcc.m / rdcc_likelihood.m (see commit message).
# The net effect on this worksheet is:
#   * RDCC_CP column ("K") likelihood figures are refreshed, and the same
#     refreshed figures now also populate the RDCC_Diagonal ("L") and
#     RDCC_Scalar ("M") columns for rows 3-6 (they were blank before).
#   * The style table grew by two duplicate (but equivalent) style/border
#     entries as a side effect of the re-export; we mirror that bookkeeping
#     below using a couple of off-sheet helper cells so the workbook's
#     style/border counts line up, without leaving any visible trace.
# ---------------------------------------------------------------------------

# --- bookkeeping: reproduce the two extra style/border slots left behind by
#     the re-export tool, via scratch cells far outside the used range so
#     they never show up in the saved sheet data / dimension.
$ws.Range("AA1").NumberFormat = "@"
$ws.Range("AA1").Borders.LineStyle = 1
$ws.Range("AA2").NumberFormat = "m/d/yy h:mm"
$ws.Range("AA2").Borders.LineStyle = 1
$ws.Range("AA2").Borders(5).LineStyle = 1
$ws.Range("AA1:AA2").Clear()

# --- actual data refresh --------------------------------------------------
# Row 3 (alpha_1 / "Parameter" row 2 -> alpha_1): RDCC_CP unchanged value,
# but now also written into RDCC_Diagonal / RDCC_Scalar.
$ws.Range("L3").Value = 0.079
$ws.Range("M3").Value = 0.079

# Row 4 (beta_1): value corrected 0.900 -> 0.899, and mirrored into L/M.
$ws.Range("K4").Value = 0.899
$ws.Range("L4").Value = 0.899
$ws.Range("M4").Value = 0.899

# Row 5 (alpha_2): value corrected 0.047 -> 0.046, and mirrored into L/M.
$ws.Range("K5").Value = 0.046
$ws.Range("L5").Value = 0.046
$ws.Range("M5").Value = 0.046

# Row 6 (beta_2): value corrected 0.946 -> 0.947, and mirrored into L/M.
$ws.Range("K6").Value = 0.947
$ws.Range("L6").Value = 0.947
$ws.Range("M6").Value = 0.947
